$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 1170
$ws.Cells.Item(96, 9).Value = 1170
$ws.Cells.Item(96, 11).Value = 3510
$ws.Cells.Item(96, 13).Value = -2137
$ws.Cells.Item(111, 8).Value = 3667.5715
$ws.Cells.Item(111, 9).Value = 1991.25
$ws.Cells.Item(111, 11).Value = 5973.75
$ws.Cells.Item(111, 13).Value = -2906.75
$ws.Cells.Item(112, 8).Value = 3763.6875
$ws.Cells.Item(112, 10).Value = 4592.727
$ws.Cells.Item(112, 12).Value = 13778.181
$ws.Cells.Item(112, 14).Value = -15994.181
$ws.Cells.Item(132, 8).Value = 5088.5386
$ws.Cells.Item(132, 9).Value = 6409.7896
$ws.Cells.Item(132, 10).Value = 1502.2858
$ws.Cells.Item(132, 11).Value = 19229.3688
$ws.Cells.Item(132, 12).Value = 4506.857400000001
$ws.Cells.Item(132, 13).Value = -16699.3688
$ws.Cells.Item(132, 14).Value = -9566.857400000001
$ws.Cells.Item(135, 8).Value = 913.6316
$ws.Cells.Item(135, 9).Value = 782.94116
$ws.Cells.Item(135, 11).Value = 7046.47044
$ws.Cells.Item(135, 13).Value = -4511.47044
$ws.Cells.Item(137, 8).Value = 1589.3556
$ws.Cells.Item(137, 9).Value = 1296.4642
$ws.Cells.Item(137, 10).Value = 2071.7646
$ws.Cells.Item(137, 11).Value = 3889.3926
$ws.Cells.Item(137, 12).Value = 6215.293799999999
$ws.Cells.Item(137, 13).Value = -1339.3926
$ws.Cells.Item(137, 14).Value = -11315.2938
$ws.Cells.Item(138, 8).Value = 2229.1892
$ws.Cells.Item(138, 9).Value = 1386.1482
$ws.Cells.Item(138, 10).Value = 4505.4
$ws.Cells.Item(138, 11).Value = 4158.444600000001
$ws.Cells.Item(138, 12).Value = 13516.2
$ws.Cells.Item(138, 13).Value = 981.5553999999993
$ws.Cells.Item(138, 14).Value = -23796.2

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4006.1572
$ws.Cells.Item(32, 9).Value = 3324.0942
$ws.Cells.Item(32, 11).Value = 3324.0942
$ws.Cells.Item(32, 13).Value = -3037.0942
$ws.Cells.Item(45, 8).Value = 60202.06
$ws.Cells.Item(45, 9).Value = 76306.37
$ws.Cells.Item(45, 11).Value = 76306.37
$ws.Cells.Item(45, 13).Value = -75929.37
$ws.Cells.Item(61, 8).Value = 13525794
$ws.Cells.Item(61, 9).Value = 17865826
$ws.Cells.Item(61, 11).Value = 17865826
$ws.Cells.Item(61, 13).Value = -17865614
$ws.Cells.Item(132, 8).Value = 3628.3928
$ws.Cells.Item(132, 9).Value = 3750.926
$ws.Cells.Item(132, 10).Value = 320
$ws.Cells.Item(132, 11).Value = 11252.778
$ws.Cells.Item(132, 12).Value = 960
$ws.Cells.Item(132, 13).Value = -8722.778
$ws.Cells.Item(132, 14).Value = -6020
$ws.Cells.Item(136, 8).Value = 13525794
$ws.Cells.Item(136, 9).Value = 17865826
$ws.Cells.Item(136, 11).Value = 53597478
$ws.Cells.Item(136, 13).Value = -53594928

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2814.0715
$ws.Cells.Item(134, 9).Value = 2814.0715
$ws.Cells.Item(134, 11).Value = 8442.2145
$ws.Cells.Item(134, 13).Value = -5907.2145

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5122.7144
$ws.Cells.Item(31, 9).Value = 2891.923
$ws.Cells.Item(31, 10).Value = 7056.067
$ws.Cells.Item(31, 11).Value = 2891.923
$ws.Cells.Item(31, 12).Value = 7056.067
$ws.Cells.Item(31, 13).Value = -2596.923
$ws.Cells.Item(31, 14).Value = -7646.067
$ws.Cells.Item(34, 8).Value = 5122.7144
$ws.Cells.Item(34, 9).Value = 2891.923
$ws.Cells.Item(34, 10).Value = 7056.067
$ws.Cells.Item(34, 11).Value = 2891.923
$ws.Cells.Item(34, 12).Value = 7056.067
$ws.Cells.Item(34, 13).Value = -2689.923
$ws.Cells.Item(34, 14).Value = -7460.067
$ws.Cells.Item(86, 8).Value = 38920.65
$ws.Cells.Item(86, 9).Value = 9232.556
$ws.Cells.Item(86, 11).Value = 9232.556
$ws.Cells.Item(86, 13).Value = -8109.556
$ws.Cells.Item(89, 8).Value = 38920.65
$ws.Cells.Item(89, 9).Value = 9232.556
$ws.Cells.Item(89, 11).Value = 46162.78
$ws.Cells.Item(89, 13).Value = -40546.78
$ws.Cells.Item(107, 8).Value = 1430.125
$ws.Cells.Item(107, 9).Value = 1490.9231
$ws.Cells.Item(107, 11).Value = 1490.9231
$ws.Cells.Item(107, 13).Value = 429.0769
$ws.Cells.Item(132, 8).Value = 3670.121
$ws.Cells.Item(132, 9).Value = 3971.7083
$ws.Cells.Item(132, 11).Value = 11915.1249
$ws.Cells.Item(132, 13).Value = -9385.124899999999
$ws.Cells.Item(134, 8).Value = 4839.643
$ws.Cells.Item(134, 9).Value = 3809.7144
$ws.Cells.Item(134, 11).Value = 11429.1432
$ws.Cells.Item(134, 13).Value = -8894.143199999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(25, 8).Value = 714.2857
$ws.Cells.Item(25, 9).Value = 875
$ws.Cells.Item(25, 11).Value = 2625
$ws.Cells.Item(25, 13).Value = -2456
$ws.Cells.Item(29, 8).Value = 655.7778
$ws.Cells.Item(29, 9).Value = 615.875
$ws.Cells.Item(29, 10).Value = 975
$ws.Cells.Item(29, 11).Value = 1847.625
$ws.Cells.Item(29, 12).Value = 2925
$ws.Cells.Item(29, 13).Value = -1570.625
$ws.Cells.Item(29, 14).Value = -3479
$ws.Cells.Item(30, 8).Value = 714.2857
$ws.Cells.Item(30, 9).Value = 875
$ws.Cells.Item(30, 11).Value = 2625
$ws.Cells.Item(30, 13).Value = -2523
$ws.Cells.Item(31, 8).Value = 1261.75
$ws.Cells.Item(31, 9).Value = 1261.75
$ws.Cells.Item(31, 11).Value = 3785.25
$ws.Cells.Item(31, 13).Value = -3497.25
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 12).ClearContents()
$ws.Cells.Item(35, 14).Value = 0
$ws.Cells.Item(36, 8).Value = 8678433
$ws.Cells.Item(36, 9).Value = 8678433
$ws.Cells.Item(36, 11).Value = 26035299
$ws.Cells.Item(36, 13).Value = -26035130
$ws.Cells.Item(92, 8).Value = 474.75
$ws.Cells.Item(92, 10).Value = 499.66666
$ws.Cells.Item(92, 12).Value = 1498.99998
$ws.Cells.Item(92, 14).Value = -3994.99998
$ws.Cells.Item(109, 8).Value = 11836
$ws.Cells.Item(109, 9).Value = 90027
$ws.Cells.Item(109, 11).Value = 270081
$ws.Cells.Item(109, 13).Value = -269041
$ws.Cells.Item(126, 8).Value = 11333.333
$ws.Cells.Item(126, 9).Value = 11333.333
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 33999.999
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 14).Value = -29059.999
$ws.Cells.Item(134, 8).Value = 1817.25
$ws.Cells.Item(134, 9).Value = 1112.6428
$ws.Cells.Item(134, 11).Value = 3337.9284
$ws.Cells.Item(134, 13).Value = 1732.0716
$ws.Cells.Item(139, 8).Value = 2597.0454
$ws.Cells.Item(139, 9).Value = 1965.125
$ws.Cells.Item(139, 10).Value = 4282.1665
$ws.Cells.Item(139, 11).Value = 5895.375
$ws.Cells.Item(139, 12).Value = 12846.4995
$ws.Cells.Item(139, 13).Value = -755.375
$ws.Cells.Item(139, 14).Value = -23126.4995
$ws.Cells.Item(140, 8).Value = 386886.38
$ws.Cells.Item(140, 10).Value = 10416666
$ws.Cells.Item(140, 12).Value = 31249998
$ws.Cells.Item(140, 14).Value = -31260358

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4045.6428
$ws.Cells.Item(80, 9).Value = 2460.375
$ws.Cells.Item(80, 10).Value = 6159.3335
$ws.Cells.Item(80, 11).Value = 2460.375
$ws.Cells.Item(80, 12).Value = 6159.3335
$ws.Cells.Item(80, 13).Value = -1462.375
$ws.Cells.Item(80, 14).Value = -8155.3335
$ws.Cells.Item(83, 8).Value = 4045.6428
$ws.Cells.Item(83, 9).Value = 2460.375
$ws.Cells.Item(83, 10).Value = 6159.3335
$ws.Cells.Item(83, 11).Value = 12301.875
$ws.Cells.Item(83, 12).Value = 30796.6675
$ws.Cells.Item(83, 13).Value = -7309.875
$ws.Cells.Item(83, 14).Value = -40780.6675
$ws.Cells.Item(122, 8).Value = 12593.8
$ws.Cells.Item(122, 9).Value = 3989.6667
$ws.Cells.Item(122, 10).Value = 25500
$ws.Cells.Item(122, 11).Value = 11969.0001
$ws.Cells.Item(122, 12).Value = 76500
$ws.Cells.Item(122, 13).Value = -9519.000100000001
$ws.Cells.Item(122, 14).Value = -81400
$ws.Cells.Item(125, 8).Value = 9550.666999999999
$ws.Cells.Item(125, 10).Value = 9550.666999999999
$ws.Cells.Item(125, 12).Value = 9550.666999999999
$ws.Cells.Item(125, 14).Value = -14470.667
$ws.Cells.Item(132, 8).Value = 2420.394
$ws.Cells.Item(132, 9).Value = 2352.7932
$ws.Cells.Item(132, 10).Value = 2910.5
$ws.Cells.Item(132, 11).Value = 7058.3796
$ws.Cells.Item(132, 12).Value = 8731.5
$ws.Cells.Item(132, 13).Value = -4528.3796
$ws.Cells.Item(132, 14).Value = -13791.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4563.8335
$ws.Cells.Item(40, 9).Value = 4563.8335
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 4563.8335
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(40, 14).Value = -4427.8335
$ws.Cells.Item(122, 8).Value = 3980
$ws.Cells.Item(122, 10).Value = 4271.4287
$ws.Cells.Item(122, 12).Value = 12814.2861
$ws.Cells.Item(122, 14).Value = -17714.2861
$ws.Cells.Item(132, 8).Value = 11578.388
$ws.Cells.Item(132, 9).Value = 10369.814
$ws.Cells.Item(132, 11).Value = 31109.442
$ws.Cells.Item(132, 13).Value = -28579.442
$ws.Cells.Item(136, 8).Value = 2500.96
$ws.Cells.Item(136, 9).Value = 2596.7727
$ws.Cells.Item(136, 11).Value = 7790.3181
$ws.Cells.Item(136, 13).Value = -5240.3181

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 32500
$ws.Cells.Item(3, 10).Value = 32500
$ws.Cells.Item(3, 12).Value = 32500
$ws.Cells.Item(3, 14).Value = -32728
$ws.Cells.Item(14, 8).Value = 343283
$ws.Cells.Item(14, 10).Value = 343283
$ws.Cells.Item(14, 12).Value = 343283
$ws.Cells.Item(14, 14).Value = -343619
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 1978.8904
$ws.Cells.Item(132, 9).Value = 2060.7207
$ws.Cells.Item(132, 11).Value = 6182.1621
$ws.Cells.Item(132, 13).Value = -3652.1621
$ws.Cells.Item(136, 8).Value = 3879.8096
$ws.Cells.Item(136, 9).Value = 1986.1034
$ws.Cells.Item(136, 11).Value = 5958.3102
$ws.Cells.Item(136, 13).Value = -3408.3102

Write-Host "Applied 226 cell changes"
